$d = $word.ActiveDocument

# Helper: force a run split at a given character offset without changing any
# text, by dropping a throw-away bookmark there and immediately removing it.
# (Word always breaks runs at bookmark boundaries; once split the runs stay
# split even after the bookmark itself is deleted.)
function Split-At($pos) {
    $r = $d.Range($pos, $pos)
    $d.Bookmarks.Add("__tmp_split__", $r)
    $d.Bookmarks("__tmp_split__").Delete()
}

# Locate the date text "DATE: 03/07/16" and figure out the offset of the
# "7" inside "03/07/16" that needs to become "6" (-> "03/06/16").
$t = $d.Content.Text
$anchor = "DATE: 03/07/16"
$idx = $t.IndexOf($anchor)

$posZero  = $idx + 9   # the "0" right before the "7" (start of old "07/" run)
$posSeven = $idx + 10  # the "7" that gets replaced with "6"
$posSlash = $idx + 11  # the "/" right after the "7"

# Replace "7" with "6".
$rng7 = $d.Range($posSeven, $posSeven + 1)
$rng7.Text = "6"

# Break the run apart so "3/" | "0" | "6" end up as separate runs, matching
# how Word splits runs around the edited text.
Split-At($posZero)
Split-At($posSeven)

# Drop the "_GoBack" bookmark (zero length) right after the newly typed "6",
# marking it as the location of the most recent edit.
$bkPos = $posSeven + 1
$rngBk = $d.Range($bkPos, $bkPos)
$d.Bookmarks.Add("_GoBack", $rngBk)

# Also split the remaining "/16" back into "/" and "16" separate runs.
Split-At($posSlash + 1)
